$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "deterministic_evo" (1st tab): the old "library" column is replaced
# by a "section" column, and the "number" column moves to the D slot with
# plain numeric values (figures are now numbered within a section instead of
# carrying a plotting-library tag).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C1").Value = "section"
$ws1.Range("D1").Value = "number"

$ws1.Range("D2").Value = 1
$ws1.Range("D3").Value = 1

$ws1.Range("D4").Select()

# ---------------------------------------------------------------------------
# Sheet "drift_langevin" (2nd tab): same header rename, plus six new rows of
# figure names for chapter 2.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Value = "section"
$ws2.Range("D1").Value = "number"

$ws2.Range("A2").Value = "buri_schematic"
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 1

$ws2.Range("A3").Value = "buri_generations"
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 2

$ws2.Range("A4").Value = "buri_experiment"
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 3

$ws2.Range("A5").Value = "ensemble_realization"
$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = 2
$ws2.Range("D5").Value = 1

$ws2.Range("A6").Value = "average_langevin"
$ws2.Range("B6").Value = 2
$ws2.Range("C6").Value = 2
$ws2.Range("D6").Value = 2

$ws2.Range("A7").Value = "random_walk"
$ws2.Range("B7").Value = 2
$ws2.Range("C7").Value = 2
$ws2.Range("D7").Value = 3

$ws2.Columns.Item(1).ColumnWidth = 17.666666666666668

$ws2.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet "classic_diffusion" (3rd tab): same header rename, plus two new rows,
# and this sheet becomes the active tab.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C1").Value = "section"
$ws3.Range("D1").Value = "number"

$ws3.Range("A2").Value = "chapman_kolmogorov"
$ws3.Range("B2").Value = 3
$ws3.Range("C2").Value = 2
$ws3.Range("D2").Value = 1

$ws3.Range("A3").Value = "schematic_reverse"
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 5
$ws3.Range("D3").Value = 1

$ws3.Columns.Item(1).ColumnWidth = 18.5

$ws3.Activate()
$ws3.Range("A4").Select()
